$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.63
$ws.Range("Q2").Value = 1.69
$ws.Range("Q3").Value = 2.58
$ws.Range("F4").Value = 2.36
$ws.Range("H4").Value = 2.78
$ws.Range("I4").Value = 2.98
$ws.Range("J4").Value = 3.9
$ws.Range("K4").Value = 4.4
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.16
$ws.Range("P4").Value = 2.74
$ws.Range("Q4").Value = 1.46
$ws.Range("R4").Value = 1.7
$ws.Range("S4").Value = 2.14
$ws.Range("U4").Value = 2.76
$ws.Range("X4").Value = 36
$ws.Range("Y4").Value = 24
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 11.5
$ws.Range("AH4").Value = 17
$ws.Range("AN4").Value = 11.5
$ws.Range("AO4").Value = 15.5
$ws.Range("F5").Value = 2.52
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 3.9
$ws.Range("J5").Value = 2.72
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1.66
$ws.Range("N5").Value = 2.2
$ws.Range("P5").Value = 1.41
$ws.Range("V5").Value = 1.36
$ws.Range("Z5").Value = 28
$ws.Range("AD5").Value = 21
$ws.Range("AE5").Value = 75
$ws.Range("AN5").Value = 70
$ws.Range("I6").Value = 8.199999999999999
$ws.Range("L6").Value = 1.28
$ws.Range("R6").Value = 1.44
$ws.Range("S6").Value = 2.84
$ws.Range("T6").Value = 1.84
$ws.Range("U6").Value = 1.94
$ws.Range("AF6").Value = 11.5
$ws.Range("F7").Value = 2.52
$ws.Range("G7").Value = 2.7
$ws.Range("H7").Value = 2.84
$ws.Range("I7").Value = 3.15
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 3.3
$ws.Range("R7").Value = 1.3
$ws.Range("S7").Value = 3.6
$ws.Range("T7").Value = 1.77
$ws.Range("U7").Value = 2.06
$ws.Range("V7").Value = 1.47
$ws.Range("W7").Value = 1.58
$ws.Range("X7").Value = 16
$ws.Range("Y7").Value = 14
$ws.Range("Z7").Value = 25
$ws.Range("AA7").Value = 65
$ws.Range("AB7").Value = 12.5
$ws.Range("AC7").Value = 9.199999999999999
$ws.Range("AD7").Value = 16
$ws.Range("AE7").Value = 44
$ws.Range("AF7").Value = 21
$ws.Range("AG7").Value = 15
$ws.Range("AH7").Value = 22
$ws.Range("AI7").Value = 60
$ws.Range("AJ7").Value = 48
$ws.Range("AK7").Value = 38
$ws.Range("AL7").Value = 55
$ws.Range("AM7").Value = 130
$ws.Range("AN7").Value = 32
$ws.Range("AO7").Value = 44
$ws.Range("F8").Value = 1.54
$ws.Range("G8").Value = 1.71
$ws.Range("I8").Value = 13.5
$ws.Range("K8").Value = 3.95
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 2.18
$ws.Range("O8").Value = 1.67
$ws.Range("P8").Value = 1.39
$ws.Range("Q8").Value = 2.78
$ws.Range("R8").Value = 1.11
$ws.Range("S8").Value = 5.3
$ws.Range("T8").Value = 2.46
$ws.Range("U8").Value = 1.3
$ws.Range("V8").Value = 1.09
$ws.Range("W8").Value = 2.4
$ws.Range("AB8").Value = 5.5
$ws.Range("G9").Value = 1.69
$ws.Range("H9").Value = 6.8
$ws.Range("I9").Value = 7.2
$ws.Range("N9").Value = 2.98
$ws.Range("O9").Value = 1.49
$ws.Range("P9").Value = 1.67
$ws.Range("U9").Value = 1.68
$ws.Range("V9").Value = 1.16
$ws.Range("W9").Value = 2.44
$ws.Range("AI9").Value = 150
$ws.Range("AO9").Value = 260
$ws.Range("F10").Value = 1.12
$ws.Range("H10").Value = 34
$ws.Range("L10").Value = 1.23
$ws.Range("P10").Value = 2.98
$ws.Range("Q10").Value = 1.46
$ws.Range("W10").Value = 8.199999999999999
$ws.Range("Z10").Value = 410
$ws.Range("AC10").Value = 990
$ws.Range("AI10").Value = 530
$ws.Range("AN10").Value = 3.05
